# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $cell = $ws.Range($CellRef)
    $origStyle = $cell.Style
    # Prefix with an apostrophe to force Excel to treat the value as text
    # (otherwise numeric-looking strings like "22.07" get coerced to numbers).
    $cell.Value = "'" + $NewValue
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '57.856.16'
Set-TextValue 'E2' '  +0.00%  '
Set-TextValue 'D3' '2.423.36'
Set-TextValue 'E3' '  -0.83%  '
Set-TextValue 'E4' '  +0.37%  '
Set-TextValue 'D5' '510.16'
Set-TextValue 'E5' '  -2.67%  '
Set-TextValue 'D6' '133.55'
Set-TextValue 'E6' '  +2.59%  '
Set-TextValue 'D7' '0.998'
Set-TextValue 'E7' '  +0.13%  '
Set-TextValue 'D8' '0.561'
Set-TextValue 'E8' '  -0.19%  '
Set-TextValue 'D9' '2.459.98'
Set-TextValue 'E9' '  +0.54%  '
Set-TextValue 'D10' '0.0986'
Set-TextValue 'E10' '  +1.31%  '
Set-TextValue 'E11' '  -1.17%  '
Set-TextValue 'D12' '0.326'
Set-TextValue 'E12' '  +1.40%  '
Set-TextValue 'E13' '  -4.82%  '
Set-TextValue 'D14' '2.861.04'
Set-TextValue 'E14' '  -0.61%  '
Set-TextValue 'D15' '57.727.82'
Set-TextValue 'E15' '  -0.01%  '
Set-TextValue 'D16' '22.07'
Set-TextValue 'E16' '  +1.99%  '
Set-TextValue 'E17' '  +1.89%  '
Set-TextValue 'D18' '2.455.01'
Set-TextValue 'E18' '  +0.62%  '
Set-TextValue 'D19' '10.40'
Set-TextValue 'E19' '  +0.63%  '
Set-TextValue 'E20' '  +0.85%  '
Set-TextValue 'D21' '316.95'
Set-TextValue 'E21' '  +1.15%  '
Set-TextValue 'D22' '6.51'
Set-TextValue 'E22' '  +6.80%  '
Set-TextValue 'D23' '0.996'
Set-TextValue 'E23' '  -0.35%  '
Set-TextValue 'E24' '  -1.95%  '
Set-TextValue 'D25' '65.49'
Set-TextValue 'E25' '  +0.92%  '
Set-TextValue 'D26' '0.995'
Set-TextValue 'E26' '  -0.63%  '
Set-TextValue 'D27' '2.519.12'
Set-TextValue 'E27' '  -1.43%  '
Set-TextValue 'E28' '  -4.93%  '
Set-TextValue 'E29' '  -1.36%  '
Set-TextValue 'D30' '7.61'
Set-TextValue 'E30' '  +5.61%  '
Set-TextValue 'D31' '173.54'
Set-TextValue 'E31' '  +0.07%  '
Set-TextValue 'D32' '0.0₃0743'
Set-TextValue 'E32' '  +1.51%  '
Set-TextValue 'D33' '1.70'
Set-TextValue 'E33' '  +0.30%  '
Set-TextValue 'D34' '6.25'
Set-TextValue 'E34' '  +2.33%  '
Set-TextValue 'D35' '1.17'
Set-TextValue 'E35' '  +2.62%  '
Set-TextValue 'E36' '  -0.16%  '
Set-TextValue 'D37' '0.989'
Set-TextValue 'E37' '  -0.90%  '
Set-TextValue 'D38' '18.17'
Set-TextValue 'E38' '  +1.96%  '
Set-TextValue 'D39' '1.28'
Set-TextValue 'E39' '  +8.72%  '
Set-TextValue 'D40' '3.88'
Set-TextValue 'E40' '  +3.66%  '
Set-TextValue 'D41' '0.819'
Set-TextValue 'E41' '  +3.56%  '
Set-TextValue 'B42' 'OKB'
Set-TextValue 'C42' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D42' '36.49'
Set-TextValue 'E42' '  +0.35%  '
Set-TextValue 'B43' 'Stacks'
Set-TextValue 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D43' '1.47'
Set-TextValue 'E43' '  +1.53%  '
Set-TextValue 'D44' '136.85'
Set-TextValue 'E44' '  +12.21%  '
Set-TextValue 'D45' '5.10'
Set-TextValue 'E45' '  +6.17%  '
Set-TextValue 'E46' '  +1.01%  '
Set-TextValue 'D47' '263.91'
Set-TextValue 'E47' '  +0.04%  '
Set-TextValue 'D48' '0.575'
Set-TextValue 'E48' '  -1.01%  '
Set-TextValue 'D49' '0.0916'
Set-TextValue 'E49' '  -0.85%  '
Set-TextValue 'D50' '0.0501'
Set-TextValue 'E50' '  +2.06%  '
Set-TextValue 'D51' '0.0216'
Set-TextValue 'E51' '  +3.05%  '

Write-Host "Applied crypto price/volume updates."
